$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="65.006.33"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +6.25%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("D3").Formula = '="2.985.17"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = '="  +4.00%  "'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("D4").Formula = '="1.00"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Formula = '="  +0.07%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("D5").Formula = '="582.17"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +2.96%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="153.99"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)

$ws.Range("E7").Formula = '="  -0.01%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)

$ws.Range("E8").Formula = '="  +2.21%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

$ws.Range("D9").Formula = '="2.981.87"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  +3.86%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

$ws.Range("E10").Formula = '="  +2.35%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$ws.Range("E11").Formula = '="  +3.32%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

$ws.Range("D12").Formula = '="0.446"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  +3.75%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

$ws.Range("D13").Formula = '="0.0000238"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +2.14%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

$ws.Range("E14").Formula = '="  +7.03%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

$ws.Range("E15").Formula = '="  +1.00%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("D16").Formula = '="65.053.40"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = '="  +6.19%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("D17").Formula = '="3.481.58"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = '="  +3.99%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

$ws.Range("D18").Formula = '="6.93"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = '="  +5.73%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("D19").Formula = '="2.986.29"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  +4.19%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

$ws.Range("D20").Formula = '="448.99"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = '="  +4.26%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)

$ws.Range("D21").Formula = '="13.68"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = '="  +5.03%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

$ws.Range("D22").Formula = '="0.679"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  +4.22%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)

$ws.Range("E23").Formula = '="  +7.33%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

$ws.Range("E24").Formula = '="  +2.75%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("D25").Formula = '="12.38"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  +5.99%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

$ws.Range("D26").Formula = '="10.74"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +7.89%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

$ws.Range("E27").Formula = '="  +10.96%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

$ws.Range("E28").Formula = '="  -0.02%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

$ws.Range("D29").Formula = '="2.43"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = '="  +18.82%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

$ws.Range("E30").Formula = '="  +12.68%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$ws.Range("E31").Formula = '="  +4.21%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

$ws.Range("E32").Formula = '="  -1.29%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

$ws.Range("E33").Formula = '="  +5.25%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)

$ws.Range("D34").Formula = '="26.74"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +5.22%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

$ws.Range("D35").Formula = '="0.999"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  -0.06%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

$ws.Range("E36").Formula = '="  +3.70%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

$ws.Range("D37").Formula = '="5.74"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  +7.06%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

$ws.Range("E38").Formula = '="  +8.45%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

$ws.Range("D39").Formula = '="45.50"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  +16.13%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

$ws.Range("D40").Formula = '="48.93"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  +0.32%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

$ws.Range("D41").Formula = '="2.91"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = '="  +3.87%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

$ws.Range("D42").Formula = '="0.302"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  +13.24%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

$ws.Range("D43").Formula = '="0.121"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  +6.67%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

$ws.Range("D44").Formula = '="8.40"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +2.44%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

$ws.Range("D45").Formula = '="389.01"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  +15.42%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)

$ws.Range("D46").Formula = '="2.764.71"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Formula = '="  +3.10%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

$ws.Range("D47").Formula = '="0.0349"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  +5.12%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

$ws.Range("D48").Formula = '="134.61"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  +0.78%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

$ws.Range("E49").Formula = '="  -0.04%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

$ws.Range("D50").Formula = '="23.26"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  +8.65%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

$ws.Range("E51").Formula = '="  +2.88%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = $false
